$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: "Height of a binary tree" (link copy/pasted from the "Delete binary tree" article)
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("A9").Value = 43908
$ws.Range("B9").Value = 43908
$ws.Range("C9").Value = "Height of a binary tree"
$ws.Range("D9").Value = "https://www.techiedelight.com/delete-given-binary-tree-iterative-recursive/"

# Row 10: "Delete binary tree"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Application.CutCopyMode = $false

$ws.Range("A10").Value = 43908
$ws.Range("B10").Value = 43908
$ws.Range("C10").Value = "Delete binary tree"
$ws.Range("D10").Value = "https://www.techiedelight.com/delete-given-binary-tree-iterative-recursive/"

$ws.Range("D10").Select()
